# "remove ovsapp in each server to constant"
#
# Each of the 4 worksheets (DC1..DC4) has a final "total" row whose A column
# contains the shared string "ovsapp" together with three numeric values in
# columns C/D/E (B is already blank). That row's contents are removed,
# leaving only the empty, still-formatted A cell behind. Once "ovsapp" is no
# longer referenced anywhere, Excel drops it from the shared-string table on
# save, which is why every other shared-string-backed cell whose index was
# greater than "ovsapp"'s shifts down by one - that happens automatically,
# we just need to clear the four "ovsapp" rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DC1")
$ws2 = $wb.Worksheets.Item("DC2")
$ws3 = $wb.Worksheets.Item("DC3")
$ws4 = $wb.Worksheets.Item("DC4")

# Clear the "ovsapp" summary row (A:E) on every sheet, keeping formatting.
$ws1.Range("A18:E18").ClearContents()
$ws2.Range("A15:E15").ClearContents()
$ws3.Range("A14:E14").ClearContents()
$ws4.Range("A18:E18").ClearContents()

# Restore each sheet's selection/view state to match the saved workbook,
# finishing on DC1 so it ends up as the active tab.
$ws2.Activate()
$ws2.Range("A15:E15").Select()
$excel.ActiveWindow.Zoom = 101

$ws3.Activate()
$ws3.Range("A14:E14").Select()

$ws4.Activate()
$ws4.Range("A18:E18").Select()

$ws1.Activate()
$ws1.Range("A18:E18").Select()
